$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (line7, line8) right after the existing "line6" row (row 7),
# pushing the "extr1".."extr8" rows down by two rows.
$ws.Rows.Item(8).Resize(2).Insert() | Out-Null

# Copy formatting (bold font + border used in column A) onto the two new rows
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A8:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Renumber the "A" (index) column for the shifted extr1..extr8 rows (now rows 10..17)
for ($i = 0; $i -le 7; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = 8 + $i
}

# Update in_service flags that changed for the shifted rows (extr1 -> row10, extr2 -> row11)
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true
